# fix bug template import
# Append a bold, red "(*)" marker to the six "required field" headers in
# row 1 of the template so users importing the sheet know which columns
# are mandatory. The six headers keep their original plain text as the
# first run and get a new run " (*)"/"(*) " in bold red Calibri 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: "Hạn xử lý" -> "Hạn xử lý " + "(*)" ---------------------------
$ws.Range("A1").Value = "Hạn xử lý (*)"
$c = $ws.Range("A1").Characters(11, 3)
$c.Font.Bold = $true
$c.Font.Color = 255

# --- B1: "Yêu cầu / Sản phẩm" -> "Yêu cầu / Sản phẩm" + " (*)" ---------
$ws.Range("B1").Value = "Yêu cầu / Sản phẩm (*)"
$c = $ws.Range("B1").Characters(19, 4)
$c.Font.Bold = $true
$c.Font.Color = 255

# --- C1: "Yêu cầu / Số lượng yêu cầu" -> "... yêu cầu " + "(*)" --------
$ws.Range("C1").Value = "Yêu cầu / Số lượng yêu cầu (*)"
$c = $ws.Range("C1").Characters(28, 3)
$c.Font.Bold = $true
$c.Font.Color = 255

# --- D1: "Yêu cầu / Đơn vị" -> "Yêu cầu / Đơn vị" + " (*)" -------------
$ws.Range("D1").Value = "Yêu cầu / Đơn vị (*)"
$c = $ws.Range("D1").Characters(17, 4)
$c.Font.Bold = $true
$c.Font.Color = 255

# --- E1: "Yêu cầu/ Từ kho" -> "Yêu cầu/ Từ kho " + "(*)" ---------------
$ws.Range("E1").Value = "Yêu cầu/ Từ kho (*)"
$c = $ws.Range("E1").Characters(17, 3)
$c.Font.Bold = $true
$c.Font.Color = 255

# --- F1: "Yêu cầu / Đến kho" -> "Yêu cầu / Đến kho " + "(*)" -----------
$ws.Range("F1").Value = "Yêu cầu / Đến kho (*)"
$c = $ws.Range("F1").Characters(19, 3)
$c.Font.Bold = $true
$c.Font.Color = 255

# The template's Page Setup was touched (orientation confirmed/written
# out explicitly as portrait) as part of the same fix.
$ws.PageSetup.Orientation = 1
